# Auto-generated Excel COM-interop script to apply cell value updates
# as described by the XML diff for Sheets/Mateus_Profits.xlsx

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1657.4286
$ws.Range("I18").Value = 2060
$ws.Range("J18").Value = 651
$ws.Range("K18").Value = 2060
$ws.Range("L18").Value = 651
$ws.Range("M18").Value = -1776
$ws.Range("N18").Value = -1219

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 6924.5
$ws.Range("I40").Value = 4255
$ws.Range("J40").Value = 7458.4
$ws.Range("K40").Value = 4255
$ws.Range("L40").Value = 7458.4
$ws.Range("M40").Value = -4080
$ws.Range("N40").Value = -7808.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 111128430
$ws.Range("I113").Value = 125003864
$ws.Range("J113").Value = 125000
$ws.Range("K113").Value = 125003864
$ws.Range("L113").Value = 125000
$ws.Range("M113").Value = -125000610
$ws.Range("N113").Value = -131508

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 4683.3335
$ws.Range("I116").Value = 3200
$ws.Range("J116").Value = 5425
$ws.Range("K116").Value = 3200
$ws.Range("L116").Value = 5425
$ws.Range("M116").Value = 242
$ws.Range("N116").Value = -12309

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 14131.889
$ws.Range("I125").Value = 13169.571
$ws.Range("J125").Value = 17500
$ws.Range("K125").Value = 118526.139
$ws.Range("L125").Value = 157500
$ws.Range("M125").Value = -116066.139
$ws.Range("N125").Value = -162420

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 2193923.2
$ws.Range("I129").Value = 767.375
$ws.Range("J129").Value = 13890755
$ws.Range("K129").Value = 2302.125
$ws.Range("L129").Value = 41672265
$ws.Range("M129").Value = 2697.875
$ws.Range("N129").Value = -41682265

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 8418.826999999999
$ws.Range("I132").Value = 1656.125
$ws.Range("J132").Value = 40879.8
$ws.Range("K132").Value = 4968.375
$ws.Range("L132").Value = 122639.4
$ws.Range("M132").Value = -2438.375
$ws.Range("N132").Value = -127699.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 4137.4585
$ws.Range("I137").Value = 3449.75
$ws.Range("J137").Value = 4275
$ws.Range("K137").Value = 10349.25
$ws.Range("L137").Value = 12825
$ws.Range("M137").Value = -7799.25
$ws.Range("N137").Value = -17925

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3810.8572
$ws.Range("I138").Value = 1433.7222
$ws.Range("J138").Value = 5191.129
$ws.Range("K138").Value = 4301.1666
$ws.Range("L138").Value = 15573.387
$ws.Range("M138").Value = 838.8334000000004
$ws.Range("N138").Value = -25853.387

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6176.4585
$ws.Range("I32").Value = 5647.5073
$ws.Range("J32").Value = 18342.334
$ws.Range("K32").Value = 5647.5073
$ws.Range("L32").Value = 18342.334
$ws.Range("M32").Value = -5360.5073
$ws.Range("N32").Value = -18916.334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3304.1333
$ws.Range("I74").Value = 2172.0264
$ws.Range("J74").Value = 9449.857
$ws.Range("K74").Value = 2172.0264
$ws.Range("L74").Value = 9449.857
$ws.Range("M74").Value = -1298.0264
$ws.Range("N74").Value = -11197.857

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 3304.1333
$ws.Range("I77").Value = 2172.0264
$ws.Range("J77").Value = 9449.857
$ws.Range("K77").Value = 10860.132
$ws.Range("L77").Value = 47249.285
$ws.Range("M77").Value = -6492.132000000001
$ws.Range("N77").Value = -55985.285

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 4982.2915
$ws.Range("I110").Value = 4454.6313
$ws.Range("J110").Value = 6987.4
$ws.Range("K110").Value = 4454.6313
$ws.Range("L110").Value = 6987.4
$ws.Range("M110").Value = -2409.6313
$ws.Range("N110").Value = -11077.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2483.3538
$ws.Range("I132").Value = 1978.6545
$ws.Range("J132").Value = 5259.2
$ws.Range("K132").Value = 5935.9635
$ws.Range("L132").Value = 15777.6
$ws.Range("M132").Value = -3405.9635
$ws.Range("N132").Value = -20837.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H134").Value = 235000
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 235000
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 235000
$ws.Range("N134").Value = -245140

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H135").Value = 82198.28999999999
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 82198.28999999999
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 82198.28999999999
$ws.Range("N135").Value = -92338.28999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1291.6552
$ws.Range("I86").Value = 1199.5454
$ws.Range("J86").Value = 1581.1428
$ws.Range("K86").Value = 1199.5454
$ws.Range("L86").Value = 1581.1428
$ws.Range("M86").Value = -76.54539999999997
$ws.Range("N86").Value = -3827.1428

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1291.6552
$ws.Range("I89").Value = 1199.5454
$ws.Range("J89").Value = 1581.1428
$ws.Range("K89").Value = 5997.727
$ws.Range("L89").Value = 7905.714
$ws.Range("M89").Value = -381.7269999999999
$ws.Range("N89").Value = -19137.714

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 14139.333
$ws.Range("I96").Value = 14139.333
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 14139.333
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -11393.333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6109.077
$ws.Range("I31").Value = 5049.5
$ws.Range("J31").Value = 6771.3125
$ws.Range("K31").Value = 5049.5
$ws.Range("L31").Value = 6771.3125
$ws.Range("M31").Value = -4754.5
$ws.Range("N31").Value = -7361.3125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 6109.077
$ws.Range("I34").Value = 5049.5
$ws.Range("J34").Value = 6771.3125
$ws.Range("K34").Value = 5049.5
$ws.Range("L34").Value = 6771.3125
$ws.Range("M34").Value = -4847.5
$ws.Range("N34").Value = -7175.3125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1795.5
$ws.Range("I94").Value = 730
$ws.Range("J94").Value = 2328.25
$ws.Range("K94").Value = 730
$ws.Range("L94").Value = 2328.25
$ws.Range("M94").Value = -279
$ws.Range("N94").Value = -3230.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("M106").ClearContents()
$ws.Range("N106").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 6339.4707
$ws.Range("I134").Value = 4537.2144
$ws.Range("J134").Value = 14750
$ws.Range("K134").Value = 13611.6432
$ws.Range("L134").Value = 44250
$ws.Range("M134").Value = -11076.6432
$ws.Range("N134").Value = -49320

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 7808700
$ws.Range("I4").Value = 7628859
$ws.Range("J4").Value = 8600000
$ws.Range("K4").Value = 22886577
$ws.Range("L4").Value = 25800000
$ws.Range("M4").Value = -22886465
$ws.Range("N4").Value = -25800224

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 5088891
$ws.Range("I7").Value = 60000000
$ws.Range("J7").Value = 96971.82000000001
$ws.Range("K7").Value = 60000000
$ws.Range("L7").Value = 96971.82000000001
$ws.Range("M7").Value = -59999888
$ws.Range("N7").Value = -97195.82000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H8").Value = 5088891
$ws.Range("I8").Value = 60000000
$ws.Range("J8").Value = 96971.82000000001
$ws.Range("K8").Value = 60000000
$ws.Range("L8").Value = 96971.82000000001
$ws.Range("M8").Value = -59999861
$ws.Range("N8").Value = -97249.82000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 651.7778
$ws.Range("I55").Value = 911.8333
$ws.Range("J55").Value = 131.66667
$ws.Range("K55").Value = 911.8333
$ws.Range("L55").Value = 131.66667
$ws.Range("M55").Value = -738.8333
$ws.Range("N55").Value = -477.66667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 6954.1304
$ws.Range("I93").Value = 2445.1516
$ws.Range("J93").Value = 18400
$ws.Range("K93").Value = 2445.1516
$ws.Range("L93").Value = 18400
$ws.Range("M93").Value = -1197.1516
$ws.Range("N93").Value = -20896

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4588.5356
$ws.Range("I136").Value = 4665.8887
$ws.Range("J136").Value = 2500
$ws.Range("K136").Value = 13997.6661
$ws.Range("L136").Value = 7500
$ws.Range("M136").Value = -11447.6661
$ws.Range("N136").Value = -12600

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 895.6
$ws.Range("I81").Value = 894.4
$ws.Range("J81").Value = 898
$ws.Range("K81").Value = 1788.8
$ws.Range("L81").Value = 1796
$ws.Range("M81").Value = -727.8
$ws.Range("N81").Value = -3918

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 895.6
$ws.Range("I84").Value = 894.4
$ws.Range("J84").Value = 898
$ws.Range("K84").Value = 8944
$ws.Range("L84").Value = 8980
$ws.Range("M84").Value = -3640
$ws.Range("N84").Value = -19588

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2815.8286
$ws.Range("I126").Value = 1773.4
$ws.Range("J126").Value = 9070.4
$ws.Range("K126").Value = 5320.200000000001
$ws.Range("L126").Value = 27211.2
$ws.Range("M126").Value = -2850.200000000001
$ws.Range("N126").Value = -32151.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3050.8538
$ws.Range("I136").Value = 2164.4595
$ws.Range("J136").Value = 11250
$ws.Range("K136").Value = 6493.3785
$ws.Range("L136").Value = 33750
$ws.Range("M136").Value = -3943.3785
$ws.Range("N136").Value = -38850
